$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row of expense data (row 8)
$ws.Range("A8").Value = "boost test boards"
$ws.Range("B8").Value = "2 units"
$ws.Range("C8").Value = 5.5
$ws.Range("D8").Value = "amazon"

# Match formatting of the row above (row 7) for columns A, B, D;
# column C should keep the $#,##0.00 style used in row 2's C cell.
$ws.Range("A8").Style = $ws.Range("A7").Style
$ws.Range("D8").Style = $ws.Range("D7").Style
$ws.Range("C8").Style = $ws.Range("C2").Style

# Update the sum formula to include the new row
$ws.Range("C10").Formula = "=SUM(C2:C8)"

# Update the selected cell
$ws.Range("B13").Select()
